$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the player names in column A (rows 2-15)
$ws.Range("A2").Value = "Vince"
$ws.Range("A3").Value = "Johnson"
$ws.Range("A4").Value = "Cena"
$ws.Range("A5").Value = "Orton"
$ws.Range("A6").Value = "Reins"
$ws.Range("A7").Value = "Rollins"
$ws.Range("A8").Value = "Streather"
$ws.Range("A9").Value = "Anrew"
$ws.Range("A10").Value = "Micheal"
$ws.Range("A11").Value = "Jordon"
$ws.Range("A12").Value = "Obama"
$ws.Range("A13").Value = "Kevin"
$ws.Range("A14").Value = "Ravi"
$ws.Range("A15").Value = "Heymen"

# Update the selected / visible range to match the author's final view state
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 6
$ws.Range("A16").Select()
